# Update generated output numbers (gh-pages data refresh) across the
# workbook's sheets, matching commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 2018
$ws1.Range("F7").Value  = 3100
$ws1.Range("F10").Value = 44
$ws1.Range("F15").Value = 10248
$ws1.Range("F21").Value = 12749
$ws1.Range("F27").Value = 604
$ws1.Range("F28").Value = 175
$ws1.Range("F34").Value = 1653
$ws1.Range("F36").Value = 72

# --- Sheet "本地生活" -------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 27

# --- Sheet "全部类型" -------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 2018
$ws4.Range("F9").Value  = 27
$ws4.Range("F10").Value = 3100
$ws4.Range("F13").Value = 44
$ws4.Range("F18").Value = 10248
$ws4.Range("F23").Value = 12749
$ws4.Range("F28").Value = 604
$ws4.Range("F30").Value = 175
$ws4.Range("F38").Value = 72
